$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (header row 1 stays unchanged):
# Row2: amaan
# Row3: p1
# Row4: p2  (replaces old rows 2,3,4 and the old row 5 "p1" record is removed)
$data = @(
    @("amaan", 15, 15, 15, 15, 15, 10, 5, 90),
    @("p1",    15, 15, 15, 15, 15, 10, 5, 90),
    @("p2",    10, 10, 10, 10, 10, 10, 2, 62)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}

# Remove the old 5th data row (worksheet row 5) so the used range becomes A1:I4
$ws.Rows.Item(5).Delete()
